$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row|Coin|Link|Volume(1h)   (columns B, C, E)
$rows = @(
    "2|Bitcoin|https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc|  +0.14%  "
    "3|Ethereum|https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth|  +0.29%  "
    "4|TetherUSD|https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt|  -0.01%  "
    "5|BNB|https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb|  +4.48%  "
    "6|USDC|https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc|  +0.00%  "
    "7|XRP|https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp|  +0.53%  "
    "8|OKB|https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb|  +1.00%  "
    "9|Cardano|https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada|  +0.62%  "
    "10|Dogecoin|https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge|  -1.56%  "
    "11|Solana|https://coinranking.com/coin/zNZHO_Sjf+solana-sol|  -1.37%  "
    "12|TRON|https://coinranking.com/coin/qUhEFk1I61atv+tron-trx|  -0.61%  "
    "13|WrappedEther|https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth|  -1.44%  "
    "14|Litecoin|https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc|  -1.84%  "
    "15|Polygon|https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic|  +3.93%  "
    "16|Polkadot|https://coinranking.com/coin/25W7FG7om+polkadot-dot|  +0.96%  "
    "17|BitcoinCash|https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch|  +2.70%  "
    "18|WrappedBTC|https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc|  -0.36%  "
    "19|Avalanche|https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax|  -3.29%  "
    "20|ShibaInu|https://coinranking.com/coin/xz24e0BjL+shibainu-shib|  -1.09%  "
    "21|Dai|https://coinranking.com/coin/MoTuySvg7+dai-dai|  -0.07%  "
    "22|WrappedliquidstakedEther2.0|https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth|  -0.70%  "
    "23|BinanceUSD|https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd|  -0.03%  "
    "24|Uniswap|https://coinranking.com/coin/_H5FVG9iW+uniswap-uni|  +0.32%  "
    "25|Chainlink|https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link|  +0.03%  "
    "26|Cosmos|https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom|  -1.97%  "
    "27|Monero|https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr|  -0.74%  "
    "28|EthereumClassic|https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc|  +0.14%  "
    "29|LidoDAOToken|https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo|  -1.12%  "
    "30|Toncoin|https://coinranking.com/coin/67YlI0K1b+toncoin-ton|  +0.85%  "
    "31|Stellar|https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm|  +0.02%  "
    "32|PancakeSwap|https://coinranking.com/coin/ncYFcP709+pancakeswap-cake|  +4.59%  "
    "33|Filecoin|https://coinranking.com/coin/ymQub4fuB+filecoin-fil|  -1.34%  "
    "34|InternetComputer(DFINITY)|https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp|  -0.30%  "
    "35|Hedera|https://coinranking.com/coin/jad286TjB+hedera-hbar|  +1.05%  "
    "36|ARBITRUM|https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb|  -0.79%  "
    "37|ImmutableX|https://coinranking.com/coin/Z96jIvLU7+immutablex-imx|  -0.54%  "
    "38|HuobiToken|https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht|  +0.11%  "
    "39|VeChain|https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet|  -1.47%  "
    "40|MXToken|https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx|  +0.14%  "
    "41|FraxShare|https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs|  +1.13%  "
    "42|Aave|https://coinranking.com/coin/ixgUfzmLR+aave-aave|  -3.70%  "
    "43|RenderToken|https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr|  -1.64%  "
    "44|TrustWalletToken|https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt|  +0.55%  "
    "45|PaxDollar|https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp|  -0.04%  "
    "46|TheSandbox|https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand|  -1.08%  "
    "47|Quant|https://coinranking.com/coin/bauj_21eYVwso+quant-qnt|  -0.72%  "
    "48|EnergySwap|https://coinranking.com/coin/SbWqqTui-+energyswap-ens|  +2.54%  "
    "49|Aptos|https://coinranking.com/coin/HGYj5JCv5+aptos-apt|  -0.16%  "
    "50|Elrond|https://coinranking.com/coin/omwkOTglq+elrond-egld|  +2.21%  "
    "51|Maker|https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr|  -5.50%  "
)

foreach ($entry in $rows) {
    $parts = $entry -split "\|", 4
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 2).Value = $parts[1]
    $ws.Cells.Item($r, 3).Value = $parts[2]
    $ws.Cells.Item($r, 5).Value = $parts[3]
}

# Column D (Price) is handled separately: values are forced to Text via
# NumberFormat "@" before assignment so numeric-looking strings (e.g.
# "1.000", "18.90", "0.7250") keep their exact textual representation
# instead of being coerced to a Double and losing trailing zeros -
# matching the source workbook, where every Price cell is inline text.
# Rows 4 and 6 (TetherUSD / USDC) keep their original "1.000" Price text
# untouched, since it does not change.
$prices = @(
    "2|30.331.31"
    "3|1.871.53"
    "5|245.15"
    "7|0.4724"
    "8|42.74"
    "9|0.2876"
    "10|0.06468"
    "11|21.09"
    "12|0.07774"
    "13|1.865.59"
    "14|95.12"
    "15|0.7250"
    "16|5.141"
    "17|275.47"
    "18|30.322.68"
    "19|13.37"
    "20|0.000007556"
    "21|0.9999"
    "22|2.108.60"
    "23|1.0000"
    "24|5.247"
    "25|6.172"
    "26|9.262"
    "27|165.59"
    "28|18.90"
    "29|1.916"
    "30|1.379"
    "31|0.09908"
    "32|1.523"
    "33|4.294"
    "34|4.032"
    "35|0.04771"
    "36|1.123"
    "37|0.6979"
    "38|2.721"
    "39|0.01845"
    "40|2.743"
    "41|6.408"
    "42|70.26"
    "43|1.916"
    "44|0.8405"
    "45|1.000"
    "46|0.4124"
    "47|102.40"
    "48|9.372"
    "49|7.099"
    "50|35.24"
    "51|919.55"
)

foreach ($entry in $prices) {
    $parts = $entry -split "\|", 2
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $parts[1]
}

Write-Output "Updated $($rows.Count) rows ($($prices.Count) price cells)"